$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update D12/E12 ---
$ws.Range("E12").Value = 0.99526005360000003
$ws.Range("D12").Formula = "=1-E12"

# --- Update D13/E13 (D13 becomes a plain value, no formula) ---
$ws.Range("D13").Value = 0.04229490617

# --- Update D14/E14 (D14 becomes a plain value, no formula) ---
$ws.Range("D14").Value = 0.04922252011

# --- Update D15/E15 (D15 becomes a plain value, no formula) ---
$ws.Range("D15").Value = 0.1757426273

# --- Clear D17 (was formula =63) ---
$ws.Range("D17").ClearContents()

# --- Update sheet view: scroll position + selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B17").Select()
